# Applies the cryptos-list price/volume refresh described in the commit.
# D-column numeric-looking text values are protected from Excel's automatic
# "smart" text-to-number coercion by briefly flipping the cell to the Text
# number format while the literal is assigned, then restoring the default
# ("Normal") cell style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

$ws.Range("D2").Value = "69.352.65"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "2.481.43"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws "D5" "565.61"
$ws.Range("E5").Value = "  -3.34%  "
Set-TextValue $ws "D6" "163.83"
$ws.Range("E6").Value = "  -5.38%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "2.478.05"
$ws.Range("E9").Value = "  -3.66%  "
Set-TextValue $ws "D10" "0.157"
$ws.Range("E10").Value = "  -6.28%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "2.934.22"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "69.239.68"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("E16").Value = "  -3.74%  "
Set-TextValue $ws "D17" "24.15"
$ws.Range("E17").Value = "  -5.64%  "
$ws.Range("D18").Value = "2.480.07"
$ws.Range("E18").Value = "  -3.62%  "
Set-TextValue $ws "D19" "11.16"
$ws.Range("E19").Value = "  -3.63%  "
Set-TextValue $ws "D20" "7.35"
$ws.Range("E20").Value = "  -7.59%  "
Set-TextValue $ws "D21" "344.49"
$ws.Range("E21").Value = "  -3.86%  "
Set-TextValue $ws "D22" "3.85"
$ws.Range("E22").Value = "  -3.10%  "
Set-TextValue $ws "D23" "1.92"
$ws.Range("E23").Value = "  -8.75%  "
$ws.Range("E24").Value = "  -0.21%  "
Set-TextValue $ws "D25" "69.71"
$ws.Range("E25").Value = "  -1.27%  "
Set-TextValue $ws "D26" "3.87"
$ws.Range("E26").Value = "  -6.57%  "
$ws.Range("D27").Value = "2.609.42"
$ws.Range("E27").Value = "  -2.87%  "
Set-TextValue $ws "D28" "8.62"
$ws.Range("E28").Value = "  -7.08%  "
Set-TextValue $ws "D29" "0.998"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0865"
$ws.Range("E30").Value = "  -7.09%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D31" "7.68"
$ws.Range("E31").Value = "  -4.16%  "
Set-TextValue $ws "D32" "441.10"
$ws.Range("E32").Value = "  -7.96%  "
$ws.Range("E33").Value = "  -8.75%  "
Set-TextValue $ws "D34" "0.999"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -4.94%  "
Set-TextValue $ws "D36" "156.64"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  -5.17%  "
Set-TextValue $ws "D38" "19.01"
$ws.Range("E38").Value = "  -0.52%  "
Set-TextValue $ws "D39" "18.06"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  -3.76%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws "D42" "4.56"
$ws.Range("E42").Value = "  -7.72%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D43" "1.57"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("E44").Value = "  -1.78%  "
Set-TextValue $ws "D45" "2.15"
$ws.Range("E45").Value = "  -10.52%  "
$ws.Range("E46").Value = "  -9.70%  "
Set-TextValue $ws "D47" "139.36"
$ws.Range("E47").Value = "  -4.68%  "
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("E49").Value = "  -6.02%  "
Set-TextValue $ws "D50" "0.0724"
$ws.Range("E50").Value = "  -2.53%  "
Set-TextValue $ws "D51" "0.570"
$ws.Range("E51").Value = "  -3.49%  "
